$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the two shared-string cable labels
$ws.Range("C19").Value = "18-1569-BLK+RED+YEL/GRN"
$ws.Range("B28").Value = "18-1569-YEL/GRN"

# Move the active selection from J8 to F14
$ws.Range("F14").Select()

# Resize columns C and D
$ws.Columns.Item(3).ColumnWidth = 32
$ws.Columns.Item(4).ColumnWidth = 15.42578125
